# Update column G ("K") values for rows 2-35 on the active sheet.
# The workbook previously stored a "Strike#" style value in column G;
# this regenerates it to hold the actual K values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 6
    3  = 5
    4  = 9
    5  = 12
    6  = 14
    7  = 7
    8  = 9
    9  = 5
    10 = 3
    11 = 16
    12 = 5
    13 = 6
    14 = 6
    15 = 12
    16 = 8
    17 = 5
    18 = 7
    19 = 10
    20 = 8
    21 = 13
    22 = 7
    23 = 7
    24 = 9
    25 = 9
    26 = 9
    27 = 13
    28 = 10
    29 = 9
    30 = 11
    31 = 9
    32 = 6
    33 = 4
    34 = 5
    35 = 2
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
